$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 47

# Column A holds a date-looking value that must be stored as plain text
# (matching the existing rows), so force text formatting, assign, then
# clear the formatting back to the sheet default to avoid leaving a
# stray style index on the new cell.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-10-01"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = 55.18000030517578
$ws.Cells.Item($row, 3).Value = 718.3499755859375
$ws.Cells.Item($row, 4).Value = 329
